$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text type on the Price/Volume columns so numeric-looking values
# (e.g. "0.9986") are written back as text, matching the source data which
# always stores these as strings (t="inlineStr"/shared string), not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns for rows 2-49
$ws.Range("D2").Value = "29.952.61"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.878.60"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "243.32"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").Value = "0.9978"
$ws.Range("D7").Value = "0.4926"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("D9").Value = "0.06602"
$ws.Range("D10").Value = "1.875.48"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "0.07194"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "0.6633"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "85.24"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "4.827"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "29.921.04"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "0.000007859"
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").Value = "0.9990"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "2.118.54"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "0.9984"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "4.748"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "9.132"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").Value = "5.549"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "147.49"
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("D26").Value = "136.68"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "16.76"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "1.910"
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "4.171"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D31").Value = "0.08625"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "3.935"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("D33").Value = "0.04964"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("D34").Value = "1.107"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").Value = "0.6991"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").Value = "2.193"
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("D38").Value = "2.678"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").Value = "0.9294"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").Value = "0.01636"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "5.980"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").Value = "0.9982"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "0.4178"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "101.86"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "7.563"
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "0.05709"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").Value = "32.49"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "8.235"
$ws.Range("E49").Value = "  -0.69%  "

# Rows 50 and 51 swap: Decentraland/Aave order flips, with new values
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "55.69"
$ws.Range("E50").Value = "  +1.59%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3695"
$ws.Range("E51").Value = "  -0.38%  "

# Strip the temporary text-number-format override so the cell style
# index matches the original (no explicit style on data cells).
$ws.Range("D2:E51").ClearFormats()
